{"js": "// Remove the duplicated 404 screenshot (+ its paragraph) and the\n// following \"Este proyecto aborda...\" introduction paragraph that sit\n// right after the \"Introducci\u00f3n\" (Heading 2) heading. The heading\n// itself is left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\n// Locate the \"Introducci\u00f3n\" heading paragraph.\nlet introHeadingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"Heading 2\" && para.text.trim() === \"Introducci\u00f3n\") {\n    introHeadingIndex = i;\n    break;\n  }\n}\n\nif (introHeadingIndex === -1) {\n  throw new Error('Could not find the \"Introducci\u00f3n\" heading paragraph.');\n}\n\n// The paragraph right after the heading holds the inline picture\n// (the broken/\"404\" screenshot), and the paragraph after that holds\n// the \"Este proyecto aborda...\" description text. Collect them first\n// (so index shifting from deleting the first one doesn't bite us),\n// then delete.\nconst pictureParagraph = paragraphs.items[introHeadingIndex + 1];\nconst textParagraph = paragraphs.items[introHeadingIndex + 2];\n\npictureParagraph.inlinePictures.load(\"items\");\ntextParagraph.load(\"text\");\nawait context.sync();\n\nif (pictureParagraph.inlinePictures.items.length === 0) {\n  throw new Error(\"Expected the paragraph after the heading to contain the screenshot picture.\");\n}\nif (!textParagraph.text.startsWith(\"Este proyecto aborda la maquetaci\u00f3n\")) {\n  throw new Error(\"Expected the paragraph after the picture to contain the intro text.\");\n}\n\ntextParagraph.delete();\npictureParagraph.delete();\nawait context.sync();\n", "ps1": "# Remove the duplicated 404 screenshot (+ its paragraph) and the\n# following \"Este proyecto aborda...\" introduction paragraph that sit\n# right after the \"Introducci\u00f3n\" (Heading 2) heading. The heading\n# itself is left untouched.\n\n$d = $word.ActiveDocument\n$paragraphs = $d.Paragraphs\n$count = $paragraphs.Count\n\n# Locate the \"Introducci\u00f3n\" heading paragraph.\n$introHeadingIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paragraphs.Item($i)\n    $styleName = $p.Style.NameLocal\n    $text = $p.Range.Text.Trim()\n    if ($styleName -eq \"Heading 2\" -and $text -eq \"Introducci\u00f3n\") {\n        $introHeadingIndex = $i\n        break\n    }\n}\n\nif ($introHeadingIndex -eq -1) {\n    throw \"Could not find the 'Introducci\u00f3n' heading paragraph.\"\n}\n\n# The paragraph right after the heading holds the inline picture (the\n# broken/\"404\" screenshot), and the paragraph after that holds the\n# \"Este proyecto aborda...\" description text.\n$picturePara = $paragraphs.Item($introHeadingIndex + 1)\n$textPara = $paragraphs.Item($introHeadingIndex + 2)\n\nif ($picturePara.Range.InlineShapes.Count -lt 1) {\n    throw \"Expected the paragraph after the heading to contain the screenshot picture.\"\n}\nif (-not $textPara.Range.Text.StartsWith(\"Este proyecto aborda\")) {\n    throw \"Expected the paragraph after the picture to contain the intro text.\"\n}\n\n# Delete the higher-indexed paragraph first so the picture paragraph's\n# index does not shift before we delete it too.\n$textPara.Range.Delete()\n$picturePara.Range.Delete()\n"}
